$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 29/30 swap: Binance-PegBSC-USD <-> RenderToken (with updated data)
$ws.Range("B29").Value = "RenderToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D29").Value = "'7.46"
$ws.Range("E29").Value = "  -3.50%  "

$ws.Range("B30").Value = "Binance-PegBSC-USD"
$ws.Range("C30").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D30").Value = "'1.00"
$ws.Range("E30").Value = "  +0.04%  "

# Price / Volume(1h) updates
$ws.Range("D2").Value = "64.801.54"
$ws.Range("E2").Value = "  -0.52%  "
$ws.Range("D3").Value = "3.511.34"
$ws.Range("E3").Value = "  -1.02%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'587.06"
$ws.Range("E5").Value = "  -1.78%  "
$ws.Range("D6").Value = "'133.34"
$ws.Range("E6").Value = "  -0.71%  "
$ws.Range("D7").Value = "3.511.02"
$ws.Range("E7").Value = "  -1.00%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  -0.70%  "
$ws.Range("E10").Value = "  +1.65%  "
$ws.Range("E11").Value = "  +1.16%  "
$ws.Range("D12").Value = "'0.386"
$ws.Range("E12").Value = "  -0.22%  "
$ws.Range("D13").Value = "4.107.66"
$ws.Range("E13").Value = "  -0.98%  "
$ws.Range("D14").Value = "'27.66"
$ws.Range("E14").Value = "  +2.78%  "
$ws.Range("D15").Value = "'0.0000181"
$ws.Range("E15").Value = "  -0.82%  "
$ws.Range("E16").Value = "  +0.76%  "
$ws.Range("D17").Value = "3.516.96"
$ws.Range("E17").Value = "  -0.98%  "
$ws.Range("D18").Value = "64.796.71"
$ws.Range("D19").Value = "'9.99"
$ws.Range("E19").Value = "  +0.78%  "
$ws.Range("D20").Value = "'14.26"
$ws.Range("E20").Value = "  -0.92%  "
$ws.Range("D21").Value = "'5.66"
$ws.Range("E21").Value = "  -2.73%  "
$ws.Range("D22").Value = "'390.76"
$ws.Range("E22").Value = "  +0.16%  "
$ws.Range("D23").Value = "'0.576"
$ws.Range("E23").Value = "  -0.21%  "
$ws.Range("D24").Value = "3.655.92"
$ws.Range("E24").Value = "  -0.97%  "
$ws.Range("D25").Value = "'74.19"
$ws.Range("E25").Value = "  +0.21%  "
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("D27").Value = "'0.0000110"
$ws.Range("E27").Value = "  -3.29%  "
$ws.Range("D28").Value = "'1.59"
$ws.Range("E28").Value = "  +2.74%  "
$ws.Range("E31").Value = "  -0.75%  "
$ws.Range("D32").Value = "'8.22"
$ws.Range("E32").Value = "  -2.82%  "
$ws.Range("D33").Value = "3.517.67"
$ws.Range("E33").Value = "  -0.94%  "
$ws.Range("E34").Value = "  +0.01%  "
$ws.Range("D35").Value = "'23.97"
$ws.Range("E35").Value = "  +0.08%  "
$ws.Range("D36").Value = "'0.146"
$ws.Range("E36").Value = "  -0.25%  "
$ws.Range("E37").Value = "  +3.74%  "
$ws.Range("D38").Value = "'171.82"
$ws.Range("E38").Value = "  +1.14%  "
$ws.Range("D39").Value = "'5.18"
$ws.Range("E39").Value = "  +3.51%  "
$ws.Range("D40").Value = "'6.96"
$ws.Range("E40").Value = "  +0.90%  "
$ws.Range("D41").Value = "'0.0808"
$ws.Range("E41").Value = "  +0.07%  "
$ws.Range("E42").Value = "  -1.13%  "
$ws.Range("D43").Value = "'26.52"
$ws.Range("E43").Value = "  +1.70%  "
$ws.Range("D44").Value = "'0.999"
$ws.Range("E44").Value = "  -0.11%  "
$ws.Range("D45").Value = "'42.26"
$ws.Range("E45").Value = "  -1.77%  "
$ws.Range("E46").Value = "  -1.16%  "
$ws.Range("D47").Value = "'4.41"
$ws.Range("E47").Value = "  -0.58%  "
$ws.Range("E48").Value = "  -0.04%  "
$ws.Range("D49").Value = "2.486.62"
$ws.Range("E49").Value = "  +1.23%  "
$ws.Range("D50").Value = "'6.87"
$ws.Range("E50").Value = "  -0.64%  "
$ws.Range("D51").Value = "'0.899"
$ws.Range("E51").Value = "  +3.04%  "
